# Auto-generated edit script: updates cryptos list (prices, % changes) and
# inserts RocketPoolETH as a new row 49, shifting RenderToken/EnergySwap down
# and dropping Aptos off the bottom of the 50-row table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.233.41'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.854.85'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7015'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07900'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3025'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08148'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '1.901.64'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.198'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7074'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '29.406.06'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.789'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007849'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").Value = '2.165.26'
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.610'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.905'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1427'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.904'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.401'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.479'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.285'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.020'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05156'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7100'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.004'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01849'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.694'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = '1.147.83'
$ws.Range("E41").Value = '  +4.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9194'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.966'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4236'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5296'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.046.86'
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.744'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.171'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.32%  '
